$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header columns (row 1)
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 2. Title-case Spanish connector words (de, del, la, las, los, el, y) in state/municipality names
$ws.Range("B6").Value = "Pabellón De Arteaga"
$ws.Range("B7").Value = "Rincón De Romos"
$ws.Range("B24").Value = "Amatenango De La Frontera"
$ws.Range("B27").Value = "Bejucal De Ocampo"
$ws.Range("B36").Value = "Comitán De Domínguez"
$ws.Range("B53").Value = "Mazapa De Madero"
$ws.Range("B55").Value = "Ocozocoautla De Espinosa"
$ws.Range("B60").Value = "San Cristóbal De Las Casas"
$ws.Range("B83").Value = "Hidalgo Del Parral"
$ws.Range("B100").Value = "San Juan De Sabinas"
$ws.Range("A110").Value = "Ciudad De México"
$ws.Range("B114").Value = "Cuajimalpa De Morelos"
$ws.Range("B128").Value = "Coneto De Comonfort"
$ws.Range("B145").Value = "San Juan Del Río"
$ws.Range("A151").Value = "Estado De México"
$ws.Range("B151").Value = "Acambay De Ruíz Castañeda"
$ws.Range("B153").Value = "Almoloya De Alquisiras"
$ws.Range("B154").Value = "Almoloya De Juárez"
$ws.Range("B166").Value = "Ecatepec De Morelos"
$ws.Range("B170").Value = "Ixtapan De La Sal"
$ws.Range("B171").Value = "Ixtapan Del Oro"
$ws.Range("B177").Value = "Naucalpan De Juárez"
$ws.Range("B184").Value = "San Felipe Del Progreso"
$ws.Range("B185").Value = "San Martín De Las Pirámides"
$ws.Range("B186").Value = "San Simón De Guerrero"
$ws.Range("B187").Value = "Soyaniquilpan De Juárez"
$ws.Range("B195").Value = "Tenango Del Aire"
$ws.Range("B196").Value = "Tenango Del Valle"
$ws.Range("B204").Value = "Tlalnepantla De Baz"
$ws.Range("B209").Value = "Villa De Allende"
$ws.Range("B216").Value = "San Miguel De Allende"
$ws.Range("B217").Value = "Apaseo El Alto"
$ws.Range("B218").Value = "Apaseo El Grande"
$ws.Range("B225").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B235").Value = "Purísima Del Rincón"
$ws.Range("B239").Value = "San Diego De La Unión"
$ws.Range("B241").Value = "San Francisco Del Rincón"
$ws.Range("B243").Value = "San Luis De La Paz"
$ws.Range("B245").Value = "Santa Cruz De Juventino Rosas"
$ws.Range("B246").Value = "Silao De La Victoria"
$ws.Range("B251").Value = "Valle De Santiago"
$ws.Range("B256").Value = "Acapulco De Juárez"
$ws.Range("B260").Value = "Atenango Del Río"
$ws.Range("B261").Value = "Atlamajalcingo Del Monte"
$ws.Range("B263").Value = "Atoyac De Álvarez"
$ws.Range("B264").Value = "Ayutla De Los Libres"
$ws.Range("B266").Value = "Buenavista De Cuéllar"
$ws.Range("B267").Value = "Chilapa De Álvarez"
$ws.Range("B268").Value = "Chilpancingo De Los Bravo"
$ws.Range("B272").Value = "Coyuca De Benítez"
$ws.Range("B273").Value = "Coyuca De Catalán"
$ws.Range("B274").Value = "Cuetzala Del Progreso"
$ws.Range("B275").Value = "Cutzamala De Pinzón"
$ws.Range("B279").Value = "Huitzuco De Los Figueroa"
$ws.Range("B280").Value = "Iguala De La Independencia"
$ws.Range("B282").Value = "Zihuatanejo De Azueta"
$ws.Range("B292").Value = "Taxco De Alarcón"
$ws.Range("B293").Value = "Técpan De Galeana"
$ws.Range("B295").Value = "Tepecoacuilco De Trujano"
$ws.Range("B297").Value = "Tixtla De Guerrero"
$ws.Range("B299").Value = "Tlapa De Comonfort"
$ws.Range("B308").Value = "Atotonilco El Grande"
$ws.Range("B312").Value = "Cuautepec De Hinojosa"
$ws.Range("B314").Value = "Huasca De Ocampo"
$ws.Range("B315").Value = "Huejutla De Reyes"
$ws.Range("B318").Value = "Jacala De Ledezma"
$ws.Range("B323").Value = "Mineral De La Reforma"
$ws.Range("B324").Value = "Mineral Del Chico"
$ws.Range("B325").Value = "Mixquiahuala De Juárez"
$ws.Range("B326").Value = "Molango De Escamilla"
$ws.Range("B328").Value = "Pachuca De Soto"
$ws.Range("B335").Value = "Tenango De Doria"
$ws.Range("B337").Value = "Tepehuacán De Guerrero"
$ws.Range("B338").Value = "Tepeji Del Río De Ocampo"
$ws.Range("B340").Value = "Tezontepec De Aldama"
$ws.Range("B344").Value = "Tula De Allende"
$ws.Range("B345").Value = "Tulancingo De Bravo"
$ws.Range("B347").Value = "Zacualtipán De Ángeles"
$ws.Range("B353").Value = "Atemajac De Brizuela"
$ws.Range("B354").Value = "Atotonilco El Alto"
$ws.Range("B355").Value = "Autlán De Navarro"
$ws.Range("B361").Value = "Concepción De Buenos Aires"
$ws.Range("B362").Value = "Cuautitlán De García Barragán"
$ws.Range("B365").Value = "Ixtlahuacán Del Río"
$ws.Range("B367").Value = "Jilotlán De Los Dolores"
$ws.Range("B370").Value = "Lagos De Moreno"
$ws.Range("B373").Value = "San Cristóbal De La Barranca"
$ws.Range("B374").Value = "San Juan De Los Lagos"
$ws.Range("B376").Value = "San Martín De Bolaños"
$ws.Range("B379").Value = "Talpa De Allende"
$ws.Range("B380").Value = "Tamazula De Gordiano"
$ws.Range("B383").Value = "Teocuitatlán De Corona"
$ws.Range("B387").Value = "Unión De San Antonio"
$ws.Range("B388").Value = "Unión De Tula"
$ws.Range("B393").Value = "Zapotlán El Grande"
$ws.Range("B449").Value = "Tiquicheo De Nicolás Romero"
$ws.Range("B475").Value = "Jonacatepec De Leandro Valle"
$ws.Range("B477").Value = "Puente De Ixtla"
$ws.Range("B481").Value = "Tetela Del Volcán"
$ws.Range("B482").Value = "Tlaltizapán De Zapata"
$ws.Range("B491").Value = "Ixtlán Del Río"
$ws.Range("B496").Value = "Santa María Del Oro"
$ws.Range("B506").Value = "Lampazos De Naranjo"
$ws.Range("B510").Value = "San Nicolás De Los Garza"
$ws.Range("B516").Value = "Constancia Del Rosario"
$ws.Range("B518").Value = "Cuyamecalco Villa De Zaragoza"
$ws.Range("B519").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B520").Value = "Heroica Ciudad De Tlaxiaco"
$ws.Range("B521").Value = "Huautla De Jiménez"
$ws.Range("B522").Value = "Ixtlán De Juárez"
$ws.Range("B523").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B528").Value = "Mariscala De Juárez"
$ws.Range("B530").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B532").Value = "Nejapa De Madero"
$ws.Range("B533").Value = "Oaxaca De Juárez"
$ws.Range("B534").Value = "Pinotepa De Don Luis"
$ws.Range("B536").Value = "Putla Villa De Guerrero"
$ws.Range("B541").Value = "San Antonino El Alto"
$ws.Range("B542").Value = "San Antonio De La Cal"
$ws.Range("B554").Value = "San José Del Progreso"
$ws.Range("B587").Value = "San Miguel Del Puerto"
$ws.Range("B596").Value = "San Pedro El Alto"
$ws.Range("B603").Value = "San Pedro Y San Pablo Ayutla"
$ws.Range("B604").Value = "San Pedro Y San Pablo Teposcolula"
$ws.Range("B619").Value = "Santa Inés Del Monte"
$ws.Range("B629").Value = "Santa María Jalapa Del Marqués"
$ws.Range("B657").Value = "Santo Domingo De Morelos"
$ws.Range("B665").Value = "Tanetze De Zaragoza"
$ws.Range("B666").Value = "Teotitlán De Flores Magón"
$ws.Range("B667").Value = "Heroica Villa Tezoatlán De Segura Y Luna, Cuna De La Independencia De Oaxaca"
$ws.Range("B668").Value = "Tlacolula De Matamoros"
$ws.Range("B669").Value = "Totontepec Villa De Morelos"
$ws.Range("B671").Value = "Villa De Tututepec"
$ws.Range("B672").Value = "Villa De Zaachila"
$ws.Range("B673").Value = "Villa Sola De Vega"
$ws.Range("B674").Value = "Villa Talea De Castro"
$ws.Range("B675").Value = "Zimatlán De Álvarez"
$ws.Range("B692").Value = "Cuayuca De Andrade"
$ws.Range("B699").Value = "Huehuetlán El Chico"
$ws.Range("B701").Value = "Ixcamilpa De Guerrero"
$ws.Range("B703").Value = "Izúcar De Matamoros"
$ws.Range("B707").Value = "Los Reyes De Juárez"
$ws.Range("B712").Value = "San Diego La Mesa Tochimiltzingo"
$ws.Range("B724").Value = "Tepanco De López"
$ws.Range("B725").Value = "Tepatlaxco De Hidalgo"
$ws.Range("B727").Value = "Tepexi De Rodríguez"
$ws.Range("B729").Value = "Tetela De Ocampo"
$ws.Range("B734").Value = "Tlacotepec De Benito Juárez"
$ws.Range("B753").Value = "Amealco De Bonfil"
$ws.Range("B755").Value = "Cadereyta De Montes"
$ws.Range("B758").Value = "Jalpan De Serra"
$ws.Range("B759").Value = "Landa De Matamoros"
$ws.Range("B760").Value = "Pinal De Amoles"
$ws.Range("B769").Value = "Armadillo De Los Infante"
$ws.Range("B774").Value = "Ciudad Del Maíz"
$ws.Range("B778").Value = "Mexquitic De Carmona"
$ws.Range("B783").Value = "San Ciro De Acosta"
$ws.Range("B785").Value = "Santa María Del Río"
$ws.Range("B786").Value = "Soledad De Graciano Sánchez"
$ws.Range("B788").Value = "Villa De Arista"
$ws.Range("B789").Value = "Villa De Arriaga"
$ws.Range("B790").Value = "Villa De Guadalupe"
$ws.Range("B791").Value = "Villa De Ramos"
$ws.Range("B792").Value = "Villa De Reyes"
$ws.Range("B839").Value = "Soto La Marina"
$ws.Range("B852").Value = "Tetla De La Solidaridad"
$ws.Range("B863").Value = "Amatlán De Los Reyes"
$ws.Range("B870").Value = "Castillo De Teayo"
$ws.Range("B880").Value = "Cosamaloapan De Carpio"
$ws.Range("B892").Value = "Hueyapan De Ocampo"
$ws.Range("B893").Value = "Ignacio De La Llave"
$ws.Range("B895").Value = "Ixhuatlán De Madero"
$ws.Range("B896").Value = "Ixhuatlán Del Café"
$ws.Range("B908").Value = "Lerdo De Tejada"
$ws.Range("B910").Value = "Martínez De La Torre"
$ws.Range("B922").Value = "Paso De Ovejas"
$ws.Range("B923").Value = "Paso Del Macho"
$ws.Range("B925").Value = "Poza Rica De Hidalgo"
$ws.Range("B931").Value = "Sayula De Alemán"
$ws.Range("B932").Value = "Soledad De Doblado"
$ws.Range("B948").Value = "Tlacotepec De Mejía"
$ws.Range("B959").Value = "Zontecomatlán De López Y Fuentes"
$ws.Range("B973").Value = "Noria De Ángeles"
$ws.Range("B977").Value = "Teúl De González Ortega"
$ws.Range("B978").Value = "Tlaltenango De Sánchez Román"
$ws.Range("B979").Value = "Villa De Cos"

# 3. Fix floating point precision (last-ULP) differences in percentage column
$ws.Range("D7").Value = 0.0009686168151879116
$ws.Range("D11").Value = 0.0009686168151879116
$ws.Range("D48").Value = 0.0009686168151879116
$ws.Range("D53").Value = 0.0009686168151879116
$ws.Range("D60").Value = 0.0009686168151879116
$ws.Range("D67").Value = 0.0009686168151879116
$ws.Range("D114").Value = 0.0009686168151879116
$ws.Range("D133").Value = 0.0009686168151879116
$ws.Range("D147").Value = 0.0009686168151879116
$ws.Range("D159").Value = 0.0009686168151879116
$ws.Range("D170").Value = 0.0009686168151879116
$ws.Range("D182").Value = 0.0009686168151879116
$ws.Range("D192").Value = 0.0009686168151879116
$ws.Range("D196").Value = 0.009492444788841536
$ws.Range("D198").Value = 0.0009686168151879116
$ws.Range("D202").Value = 0.0009686168151879116
$ws.Range("D225").Value = 0.0009686168151879116
$ws.Range("D262").Value = 0.0009686168151879116
$ws.Range("D275").Value = 0.0009686168151879116
$ws.Range("D310").Value = 0.0009686168151879116
$ws.Range("D336").Value = 0.0009686168151879116
$ws.Range("D345").Value = 0.0009686168151879116
$ws.Range("D360").Value = 0.0009686168151879116
$ws.Range("D369").Value = 0.0009686168151879116
$ws.Range("D390").Value = 0.0009686168151879116
$ws.Range("D421").Value = 0.0009686168151879116
$ws.Range("D427").Value = 0.0009686168151879116
$ws.Range("D454").Value = 0.0009686168151879116
$ws.Range("D481").Value = 0.0009686168151879116
$ws.Range("D514").Value = 0.0009686168151879116
$ws.Range("D515").Value = 0.0009686168151879116
$ws.Range("D518").Value = 0.0009686168151879116
$ws.Range("D523").Value = 0.0009686168151879116
$ws.Range("D536").Value = 0.0009686168151879116
$ws.Range("D569").Value = 0.0009686168151879116
$ws.Range("D600").Value = 0.0009686168151879116
$ws.Range("D609").Value = 0.0009686168151879116
$ws.Range("D617").Value = 0.0009686168151879116
$ws.Range("D664").Value = 0.0009686168151879116
$ws.Range("D669").Value = 0.0009686168151879116
$ws.Range("D675").Value = 0.009492444788841536
$ws.Range("D687").Value = 0.0009686168151879116
$ws.Range("D755").Value = 0.0009686168151879116
$ws.Range("D759").Value = 0.0009686168151879116
$ws.Range("D764").Value = 0.009492444788841536
$ws.Range("D771").Value = 0.0009686168151879116
$ws.Range("D774").Value = 0.0009686168151879116
$ws.Range("D796").Value = 0.0009686168151879116
$ws.Range("D814").Value = 0.0009686168151879116
$ws.Range("D819").Value = 0.0009686168151879116
$ws.Range("D827").Value = 0.0009686168151879116
$ws.Range("D836").Value = 0.0009686168151879116
$ws.Range("D871").Value = 0.0009686168151879116
$ws.Range("D906").Value = 0.0009686168151879116
$ws.Range("D910").Value = 0.0009686168151879116
$ws.Range("D923").Value = 0.0009686168151879116
$ws.Range("D932").Value = 0.0009686168151879116
$ws.Range("D958").Value = 0.0009686168151879116

# 4. Remove footer/metadata rows 984-989 (blank row + notes), shrinking the used range to A1:D983
$ws.Range("A984:D989").EntireRow.Delete()

Write-Host "Edit complete"